# Applies the edit described by the commit "Add files via upload":
#   1. On the "Oversikt" slide, add a new first bullet "Why automate?"
#      before the existing "High level overview" bullet.
#   2. Insert a brand-new "Q&A" / "Questions" slide right after the
#      "Tips" slide (and before the "Extras" slide).

$p = $ppt.ActivePresentation

# --- 1) "Oversikt" slide: add "Why automate?" as the new first bullet ---
$oversiktSlide = $p.Slides.Item(2)
$contentRange = $oversiktSlide.Shapes.Item(2).TextFrame.TextRange
[void]$contentRange.InsertBefore("Why automate?`r")

# --- 2) Insert new "Q&A" slide at position 14 (after "Tips", before "Extras") ---
$newSlide = $p.Slides.Add(14, 2)

$titleRange = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "Q&A"
$titleRange.LanguageID = "nb-NO"

$bodyRange = $newSlide.Shapes.Item(2).TextFrame.TextRange
$bodyRange.Text = "Questions"
$bodyRange.LanguageID = "nb-NO"
